# edit.ps1
# Applies the diff: restructure intro/target heading, move lastRenderedPageBreak,
# remove stray bookmark, and append a new "About" (CV) section at the end.

$d = $word.ActiveDocument

function Wrap-Xml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Replace the last paragraph (drop its trailing bookmarkEnd) and append the
#    brand new "About" / CV section right after it.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastPara.Range.InsertXML((Wrap-Xml '<w:p><w:pPr><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>The way the score is handled is using the same metrics as previously used. But instead of it constantly showing a score, it instead awards points for good behaviour</w:t></w:r><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> and actions whilst driving, and subtracts points when bad behaviour</w:t></w:r><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> or actions</w:t></w:r><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> are done. This results in a score per distance. Which the user can see and improve upon. Secondly the system also provides feedback as to why it added or subtracted points such as taking gentle corners or harsh braking.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kop1"/><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>About</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">To me nothing is more satisfying than sharing our thoughts and creations with others. Especially if at the end the other person walks off with happy thoughts and a smile. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">I am currently in the final year of the User Experience Design program at The Hague University of Applied Sciences. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Where I put emphasis on understanding users. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Working as a freelancer alongside to improve myself in various fields. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Born 21 years ago in a village called Zwijndrecht on the outskirts of the Dutch Randstad. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Early in my youth, my parents and me moved to the neighbouring Hendrik-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Ido</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">-Ambacht, where I still live to this day. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">During my high school period I obtained a CAE certificate and an International Baccalaureate, which allow me to operate on a higher level of English. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">A Technasium certificate, which gave me an early taste at designing. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">I finished high school with VWO NG&amp;NT. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>When I’m not busy with work I can be found playing video-games or watching series and movies. But mostly I can be found reading up on what interest me. Be it work-related, game-related or a documentary on space.</w:t></w:r></w:p>'))

# ---------------------------------------------------------------------------
# 2) Drop the bookmarkStart that used to precede "The final redesign..." run
#    (the _GoBack bookmark now lives at the top of the document, around the
#    new "Target" heading).
# ---------------------------------------------------------------------------
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "The final redesign of the TomTom Curfer*") {
        $targetParaIndex = $i
        break
    }
}
$p49 = $d.Paragraphs($targetParaIndex)
$p49.Range.InsertXML((Wrap-Xml '<w:p><w:pPr><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>The final redesign of the TomTom Curfer consists of a home screen where your driving score is determined on you average driving score divided by your total distance travelled. A group section where people can join and leave groups and a score overview page.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>br</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>&gt;</w:t></w:r></w:p>'))

# ---------------------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from "The project" heading to "The
#    brief" heading, within the Secret Motion section (2nd occurrence of
#    each heading in the document).
# ---------------------------------------------------------------------------
$briefParaIndex = -1
$briefCount = 0
$projectParaIndex = -1
$projectCount = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -eq "The brief`r") {
        $briefCount = $briefCount + 1
        if ($briefCount -eq 2) { $briefParaIndex = $i }
    }
    if ($t -eq "The project`r") {
        $projectCount = $projectCount + 1
        if ($projectCount -eq 2) { $projectParaIndex = $i }
    }
}

$pProject = $d.Paragraphs($projectParaIndex)
$pProject.Range.InsertXML((Wrap-Xml '<w:p><w:pPr><w:pStyle w:val="Kop2"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>The project</w:t></w:r></w:p>'))

$pBrief = $d.Paragraphs($briefParaIndex)
$pBrief.Range.InsertXML((Wrap-Xml '<w:p><w:pPr><w:pStyle w:val="Kop2"/><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t>The brief</w:t></w:r></w:p>'))

# ---------------------------------------------------------------------------
# 4) Split the opening "Target" heading paragraph into: a new "Intro"
#    heading, an intro paragraph about being a UX designer, and a fresh
#    "Target" heading that now also carries the _GoBack bookmark.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML((Wrap-Xml '<w:p><w:pPr><w:pStyle w:val="Kop1"/><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>Intro</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>I am a Dutch user experience designer, focusing on putting a smile on the face of the user. Understanding users to make amazing designs. That is what I do.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kop1"/><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-GB"/></w:rPr><w:t>Target</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'))

Write-Host "done"
